# Auto-generated: update leve-flip profit figures (scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5412.875
$ws.Range("I76").Value = 4930.6
$ws.Range("K76").Value = 4930.6
$ws.Range("M76").Value = -4615.6

$ws.Range("H79").Value = 5412.875
$ws.Range("I79").Value = 4930.6
$ws.Range("K79").Value = 4930.6
$ws.Range("M79").Value = -3838.6

$ws.Range("H135").Value = 1324.5333
$ws.Range("I135").Value = 613.7273
$ws.Range("K135").Value = 5523.545700000001
$ws.Range("M135").Value = -2988.545700000001

$ws.Range("H137").Value = 72524.875
$ws.Range("I137").Value = 185766.5
$ws.Range("J137").Value = 4579.9
$ws.Range("K137").Value = 557299.5
$ws.Range("L137").Value = 13739.7
$ws.Range("M137").Value = -554749.5
$ws.Range("N137").Value = -18839.7

$ws.Range("H138").Value = 3078.5977
$ws.Range("J138").Value = 3552.9219
$ws.Range("L138").Value = 10658.7657
$ws.Range("N138").Value = -20938.7657


# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1768.3334
$ws.Range("I45").Value = 1536.6666
$ws.Range("K45").Value = 1536.6666
$ws.Range("M45").Value = -1159.6666

$ws.Range("H61").Value = 3851.8718
$ws.Range("I61").Value = 2239.0322
$ws.Range("K61").Value = 2239.0322
$ws.Range("M61").Value = -2027.0322

$ws.Range("H74").Value = 51406.4
$ws.Range("I74").Value = 71878.78999999999
$ws.Range("K74").Value = 71878.78999999999
$ws.Range("M74").Value = -71004.78999999999

$ws.Range("H77").Value = 51406.4
$ws.Range("I77").Value = 71878.78999999999
$ws.Range("K77").Value = 359393.95
$ws.Range("M77").Value = -355025.95

$ws.Range("H132").Value = 2011.5
$ws.Range("I132").Value = 2095.139
$ws.Range("K132").Value = 6285.417
$ws.Range("M132").Value = -3755.417

$ws.Range("H136").Value = 3851.8718
$ws.Range("I136").Value = 2239.0322
$ws.Range("K136").Value = 6717.096600000001
$ws.Range("M136").Value = -4167.096600000001


# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4088.7222
$ws.Range("I99").Value = 4088.7222
$ws.Range("K99").Value = 4088.7222
$ws.Range("M99").Value = -2590.7222

$ws.Range("H107").Value = 1410.75
$ws.Range("I107").Value = 1548.8334
$ws.Range("K107").Value = 1548.8334
$ws.Range("M107").Value = 371.1666


# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1837.7826
$ws.Range("I58").Value = 2247.6155
$ws.Range("J58").Value = 1305
$ws.Range("K58").Value = 2247.6155
$ws.Range("L58").Value = 1305
$ws.Range("M58").Value = -2044.6155
$ws.Range("N58").Value = -1711

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H109").Value = 45555
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H122").Value = 7999.75
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H136").Value = 1837.7826
$ws.Range("I136").Value = 2247.6155
$ws.Range("J136").Value = 1305
$ws.Range("K136").Value = 6742.8465
$ws.Range("L136").Value = 3915
$ws.Range("M136").Value = -4192.8465
$ws.Range("N136").Value = -9015


# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 684.93335
$ws.Range("I114").Value = 678.7778
$ws.Range("K114").Value = 2036.3334
$ws.Range("M114").Value = 1217.6666

$ws.Range("H117").Value = 623.3
$ws.Range("I117").Value = 526.8
$ws.Range("J117").Value = 719.8
$ws.Range("K117").Value = 1580.4
$ws.Range("L117").Value = 2159.4
$ws.Range("M117").Value = 1861.6
$ws.Range("N117").Value = -9043.4

$ws.Range("H121").Value = 52636840
$ws.Range("J121").Value = 7539.5
$ws.Range("L121").Value = 22618.5
$ws.Range("N121").Value = -25238.5

$ws.Range("H128").Value = 285810.9
$ws.Range("I128").Value = 285810.9
$ws.Range("K128").Value = 857432.7000000001
$ws.Range("M128").Value = -852452.7000000001

$ws.Range("H137").Value = 4223.75
$ws.Range("J137").Value = 4998.3335
$ws.Range("L137").Value = 14995.0005
$ws.Range("N137").Value = -25195.0005


# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3156.4583
$ws.Range("I102").Value = 3121.4211
$ws.Range("K102").Value = 3121.4211
$ws.Range("M102").Value = -1499.4211

$ws.Range("H122").Value = 4687.615
$ws.Range("I122").Value = 4538.636
$ws.Range("K122").Value = 13615.908
$ws.Range("M122").Value = -11165.908

$ws.Range("H132").Value = 92682
$ws.Range("I132").Value = 144215
$ws.Range("J132").Value = 2499.25
$ws.Range("K132").Value = 432645
$ws.Range("L132").Value = 7497.75
$ws.Range("M132").Value = -430115
$ws.Range("N132").Value = -12557.75


# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 45699.2
$ws.Range("J38").Value = 45249.75
$ws.Range("L38").Value = 45249.75
$ws.Range("N38").Value = -46069.75

$ws.Range("H43").Value = 18059.611
$ws.Range("J43").Value = 17671.77
$ws.Range("L43").Value = 17671.77
$ws.Range("N43").Value = -18057.77

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()


# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 17200.8
$ws.Range("I122").Value = 3396.5
$ws.Range("K122").Value = 10189.5
$ws.Range("M122").Value = -7739.5

$ws.Range("H132").Value = 1433.56
$ws.Range("I132").Value = 1389.125
$ws.Range("K132").Value = 4167.375
$ws.Range("M132").Value = -1637.375

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 210944.5
$ws.Range("I136").Value = 210944.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 632833.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -630283.5
$ws.Range("N136").ClearContents()

$ws.Range("H137").Value = 3400
$ws.Range("I137").Value = 3400
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 3400
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 1700
$ws.Range("N137").ClearContents()

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

